$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "cluster2"
$ws.Range("B7").Value = "NetApp"
$ws.Range("C7").Value = "svm1_cluster2"
